# Applies the cryptos-list refresh described in the commit message.
# Price (column D) and Volume(1h) (column E) values are updated per row;
# a handful of rows (36/37, 41/42, 50/51) were additionally re-ranked, so
# their Coin name (B) and Link (C) swapped along with D/E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.052.31"
$ws.Range("E2").Value = "  -0.36%  "

# Row 3
$ws.Range("D3").Value = "2.052.28"
$ws.Range("E3").Value = "  -0.33%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "'246.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "

# Row 6
$ws.Range("D6").Value = "'0.659"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.64%  "

# Row 7
$ws.Range("D7").Value = "'58.40"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.07%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  -2.34%  "

# Row 10
$ws.Range("D10").Value = "'0.0773"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.63%  "

# Row 11
$ws.Range("D11").Value = "'0.110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.06%  "

# Row 12
$ws.Range("D12").Value = "'15.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.66%  "

# Row 13
$ws.Range("D13").Value = "'0.893"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.70%  "

# Row 14
$ws.Range("D14").Value = "2.349.83"
$ws.Range("E14").Value = "  -0.28%  "

# Row 16
$ws.Range("D16").Value = "2.073.44"
$ws.Range("E16").Value = "  +0.81%  "

# Row 17
$ws.Range("D17").Value = "'18.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.34%  "

# Row 18
$ws.Range("D18").Value = "37.037.98"
$ws.Range("E18").Value = "  -0.36%  "

# Row 19
$ws.Range("D19").Value = "'73.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.31%  "

# Row 20
$ws.Range("E20").Value = "  -2.32%  "

# Row 21
$ws.Range("D21").Value = "'5.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.81%  "

# Row 22
$ws.Range("D22").Value = "'238.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "

# Row 23
$ws.Range("E23").Value = "  +0.01%  "

# Row 24
$ws.Range("E24").Value = "  +0.92%  "

# Row 25
$ws.Range("D25").Value = "'9.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.81%  "

# Row 26
$ws.Range("D26").Value = "'168.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "

# Row 27
$ws.Range("E27").Value = "  -3.47%  "

# Row 28
$ws.Range("D28").Value = "'20.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.47%  "

# Row 29
$ws.Range("E29").Value = "  +15.32%  "

# Row 30
$ws.Range("E30").Value = "  -1.03%  "

# Row 31
$ws.Range("E31").Value = "  -2.08%  "

# Row 32
$ws.Range("D32").Value = "'4.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.44%  "

# Row 33
$ws.Range("E33").Value = "  -1.49%  "

# Row 34
$ws.Range("E34").Value = "  +0.10%  "

# Row 35
$ws.Range("D35").Value = "'1.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.26%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.0849"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.52%  "

# Row 37
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.25%  "

# Row 38
$ws.Range("E38").Value = "  -3.56%  "

# Row 39
$ws.Range("E39").Value = "  -0.95%  "

# Row 40
$ws.Range("D40").Value = "'3.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.47%  "

# Row 41
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0223"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.97%  "

# Row 42
$ws.Range("B42").Value = "Cronos"
$ws.Range("C42").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D42").Value = "'0.0975"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -11.14%  "

# Row 43
$ws.Range("E43").Value = "  +0.58%  "

# Row 44
$ws.Range("D44").Value = "'97.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.34%  "

# Row 45
$ws.Range("E45").Value = "  -6.74%  "

# Row 46
$ws.Range("D46").Value = "1.301.23"
$ws.Range("E46").Value = "  +0.23%  "

# Row 47
$ws.Range("E47").Value = "  -5.23%  "

# Row 48
$ws.Range("D48").Value = "'2.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
$ws.Range("E49").Value = "  -1.71%  "

# Row 50
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.230.89"
$ws.Range("E50").Value = "  -0.69%  "

# Row 51
$ws.Range("B51").Value = "FTXToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D51").Value = "'3.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.66%  "
